# Auto-generated edit script: updates computed profit/price columns (H-N)
# for specific Leve rows across multiple sheets, per the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value2 = 2900035.5
$ws.Range("I70").Value2 = 1602
$ws.Range("J70").Value2 = 3031782.5
$ws.Range("K70").Value2 = 4806
$ws.Range("L70").Value2 = 9095347.5
$ws.Range("M70").Value2 = -4536
$ws.Range("N70").Value2 = -9095887.5
$ws.Range("H73").Value2 = 2900035.5
$ws.Range("I73").Value2 = 1602
$ws.Range("J73").Value2 = 3031782.5
$ws.Range("K73").Value2 = 4806
$ws.Range("L73").Value2 = 9095347.5
$ws.Range("M73").Value2 = -3870
$ws.Range("N73").Value2 = -9097219.5
$ws.Range("H98").Value2 = 662.5
$ws.Range("I98").Value2 = 639.0625
$ws.Range("J98").Value2 = 725
$ws.Range("K98").Value2 = 639.0625
$ws.Range("L98").Value2 = 725
$ws.Range("M98").Value2 = 858.9375
$ws.Range("N98").Value2 = -3721
$ws.Range("H107").Value2 = 4376.6665
$ws.Range("I107").Value2 = 1032.56
$ws.Range("J107").Value2 = 21097.2
$ws.Range("K107").Value2 = 1032.56
$ws.Range("L107").Value2 = 21097.2
$ws.Range("M107").Value2 = 887.4400000000001
$ws.Range("N107").Value2 = -24937.2
$ws.Range("H116").Value2 = 3500107.8
$ws.Range("I116").Value2 = 4811717
$ws.Range("J116").Value2 = 2483.3333
$ws.Range("K116").Value2 = 4811717
$ws.Range("L116").Value2 = 2483.3333
$ws.Range("M116").Value2 = -4808275
$ws.Range("N116").Value2 = -9367.3333
$ws.Range("H122").Value2 = 662.5
$ws.Range("I122").Value2 = 639.0625
$ws.Range("J122").Value2 = 725
$ws.Range("K122").Value2 = 1917.1875
$ws.Range("L122").Value2 = 2175
$ws.Range("M122").Value2 = 532.8125
$ws.Range("N122").Value2 = -7075
$ws.Range("H137").Value2 = 36577.1
$ws.Range("I137").Value2 = 1326.5
$ws.Range("K137").Value2 = 3979.5
$ws.Range("M137").Value2 = -1429.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 1146.7142
$ws.Range("I2").Value2 = 1185.9166
$ws.Range("J2").Value2 = 911.5
$ws.Range("K2").Value2 = 1185.9166
$ws.Range("L2").Value2 = 911.5
$ws.Range("M2").Value2 = -1072.9166
$ws.Range("N2").Value2 = -1137.5
$ws.Range("H32").Value2 = 11242121
$ws.Range("I32").Value2 = 14494384
$ws.Range("J32").Value2 = 21816.25
$ws.Range("K32").Value2 = 14494384
$ws.Range("L32").Value2 = 21816.25
$ws.Range("M32").Value2 = -14494097
$ws.Range("N32").Value2 = -22390.25
$ws.Range("H45").Value2 = 2461.923
$ws.Range("I45").Value2 = 2335.158
$ws.Range("J45").Value2 = 2806
$ws.Range("K45").Value2 = 2335.158
$ws.Range("L45").Value2 = 2806
$ws.Range("M45").Value2 = -1958.158
$ws.Range("N45").Value2 = -3560
$ws.Range("H61").Value2 = 5780.9585
$ws.Range("I61").Value2 = 6361.45
$ws.Range("J61").Value2 = 2878.5
$ws.Range("K61").Value2 = 6361.45
$ws.Range("L61").Value2 = 2878.5
$ws.Range("M61").Value2 = -6149.45
$ws.Range("N61").Value2 = -3302.5
$ws.Range("H116").Value2 = 1146.7142
$ws.Range("I116").Value2 = 1185.9166
$ws.Range("J116").Value2 = 911.5
$ws.Range("K116").Value2 = 1185.9166
$ws.Range("L116").Value2 = 911.5
$ws.Range("M116").Value2 = 1108.0834
$ws.Range("N116").Value2 = -5499.5
$ws.Range("H132").Value2 = 2125.9678
$ws.Range("I132").Value2 = 1848.4445
$ws.Range("J132").Value2 = 3999.25
$ws.Range("K132").Value2 = 5545.333500000001
$ws.Range("L132").Value2 = 11997.75
$ws.Range("M132").Value2 = -3015.333500000001
$ws.Range("N132").Value2 = -17057.75
$ws.Range("H136").Value2 = 5780.9585
$ws.Range("I136").Value2 = 6361.45
$ws.Range("J136").Value2 = 2878.5
$ws.Range("K136").Value2 = 19084.35
$ws.Range("L136").Value2 = 8635.5
$ws.Range("M136").Value2 = -16534.35
$ws.Range("N136").Value2 = -13735.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 1146.7142
$ws.Range("I3").Value2 = 1185.9166
$ws.Range("J3").Value2 = 911.5
$ws.Range("K3").Value2 = 1185.9166
$ws.Range("L3").Value2 = 911.5
$ws.Range("M3").Value2 = -1071.9166
$ws.Range("N3").Value2 = -1139.5
$ws.Range("H80").Value2 = 263.5
$ws.Range("I80").Value2 = 208.5
$ws.Range("J80").Value2 = 296.5
$ws.Range("K80").Value2 = 208.5
$ws.Range("L80").Value2 = 296.5
$ws.Range("M80").Value2 = 789.5
$ws.Range("N80").Value2 = -2292.5
$ws.Range("H83").Value2 = 263.5
$ws.Range("I83").Value2 = 208.5
$ws.Range("J83").Value2 = 296.5
$ws.Range("K83").Value2 = 1042.5
$ws.Range("L83").Value2 = 1482.5
$ws.Range("M83").Value2 = 3949.5
$ws.Range("N83").Value2 = -11466.5
$ws.Range("H107").Value2 = 2242.4285
$ws.Range("I107").Value2 = 2329.3076
$ws.Range("J107").Value2 = 1113
$ws.Range("K107").Value2 = 2329.3076
$ws.Range("L107").Value2 = 1113
$ws.Range("M107").Value2 = -409.3076000000001
$ws.Range("N107").Value2 = -4953
$ws.Range("H134").Value2 = 2083.353
$ws.Range("I134").Value2 = 2185.3572
$ws.Range("J134").Value2 = 1607.3334
$ws.Range("K134").Value2 = 6556.071599999999
$ws.Range("L134").Value2 = 4822.0002
$ws.Range("M134").Value2 = -4021.071599999999
$ws.Range("N134").Value2 = -9892.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 1287.5306
$ws.Range("I31").Value2 = 1173.3704
$ws.Range("J31").Value2 = 1427.6364
$ws.Range("K31").Value2 = 1173.3704
$ws.Range("L31").Value2 = 1427.6364
$ws.Range("M31").Value2 = -878.3704
$ws.Range("N31").Value2 = -2017.6364
$ws.Range("H34").Value2 = 1287.5306
$ws.Range("I34").Value2 = 1173.3704
$ws.Range("J34").Value2 = 1427.6364
$ws.Range("K34").Value2 = 1173.3704
$ws.Range("L34").Value2 = 1427.6364
$ws.Range("M34").Value2 = -971.3704
$ws.Range("N34").Value2 = -1831.6364
$ws.Range("H58").Value2 = 2226.5881
$ws.Range("I58").Value2 = 2401.4285
$ws.Range("J58").Value2 = 2104.2
$ws.Range("K58").Value2 = 2401.4285
$ws.Range("L58").Value2 = 2104.2
$ws.Range("M58").Value2 = -2198.4285
$ws.Range("N58").Value2 = -2510.2
$ws.Range("H107").Value2 = 556444.6
$ws.Range("I107").Value2 = 1250781.9
$ws.Range("J107").Value2 = 974.8
$ws.Range("K107").Value2 = 1250781.9
$ws.Range("L107").Value2 = 974.8
$ws.Range("M107").Value2 = -1248861.9
$ws.Range("N107").Value2 = -4814.8
$ws.Range("H136").Value2 = 2226.5881
$ws.Range("I136").Value2 = 2401.4285
$ws.Range("J136").Value2 = 2104.2
$ws.Range("K136").Value2 = 7204.2855
$ws.Range("L136").Value2 = 6312.599999999999
$ws.Range("M136").Value2 = -4654.2855
$ws.Range("N136").Value2 = -11412.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 611.43335
$ws.Range("I5").Value2 = 716.0909
$ws.Range("J5").Value2 = 550.8421
$ws.Range("K5").Value2 = 2148.2727
$ws.Range("L5").Value2 = 1652.5263
$ws.Range("M5").Value2 = -2036.2727
$ws.Range("N5").Value2 = -1876.5263
$ws.Range("H68").Value2 = 906.7646999999999
$ws.Range("I68").Value2 = 764.97144
$ws.Range("J68").Value2 = 1057.1515
$ws.Range("K68").Value2 = 2294.91432
$ws.Range("L68").Value2 = 3171.4545
$ws.Range("M68").Value2 = -1483.91432
$ws.Range("N68").Value2 = -4793.4545
$ws.Range("H71").Value2 = 906.7646999999999
$ws.Range("I71").Value2 = 764.97144
$ws.Range("J71").Value2 = 1057.1515
$ws.Range("K71").Value2 = 6884.74296
$ws.Range("L71").Value2 = 9514.363499999999
$ws.Range("M71").Value2 = -2828.74296
$ws.Range("N71").Value2 = -17626.3635
$ws.Range("H135").Value2 = 611.43335
$ws.Range("I135").Value2 = 716.0909
$ws.Range("J135").Value2 = 550.8421
$ws.Range("K135").Value2 = 6444.8181
$ws.Range("L135").Value2 = 4957.5789
$ws.Range("M135").Value2 = -3909.8181
$ws.Range("N135").Value2 = -10027.5789

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value2 = 3440.9285
$ws.Range("I122").Value2 = 3882.8
$ws.Range("J122").Value2 = 3195.4443
$ws.Range("K122").Value2 = 11648.4
$ws.Range("L122").Value2 = 9586.332900000001
$ws.Range("M122").Value2 = -9198.400000000001
$ws.Range("N122").Value2 = -14486.3329
$ws.Range("H132").Value2 = 4767.646
$ws.Range("I132").Value2 = 4849.6587
$ws.Range("J132").Value2 = 4287.2856
$ws.Range("K132").Value2 = 14548.9761
$ws.Range("L132").Value2 = 12861.8568
$ws.Range("M132").Value2 = -12018.9761
$ws.Range("N132").Value2 = -17921.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value2 = 2645.7273
$ws.Range("I82").Value2 = 2325.375
$ws.Range("J82").Value2 = 3500
$ws.Range("K82").Value2 = 2325.375
$ws.Range("L82").Value2 = 3500
$ws.Range("M82").Value2 = -1964.375
$ws.Range("N82").Value2 = -4222
$ws.Range("H85").Value2 = 2645.7273
$ws.Range("I85").Value2 = 2325.375
$ws.Range("J85").Value2 = 3500
$ws.Range("K85").Value2 = 2325.375
$ws.Range("L85").Value2 = 3500
$ws.Range("M85").Value2 = -1077.375
$ws.Range("N85").Value2 = -5996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 1385.3636
$ws.Range("I132").Value2 = 1088.8077
$ws.Range("J132").Value2 = 2486.8572
$ws.Range("K132").Value2 = 3266.4231
$ws.Range("L132").Value2 = 7460.571599999999
$ws.Range("M132").Value2 = -736.4231
$ws.Range("N132").Value2 = -12520.5716

